# Apply the LinuxForHealth rebrand / version bump edit to the
# StructureDefinition-based-on-value workbook.
#
# Sheet "Metadata": URL / Version / Date / Publisher property values.
# Sheet "Elements": the "based-on-value" URL repeated in the fixed-value
#   column for Extension.url, plus the FHIR ele-1/ext-1 constraint text
#   moving from the top-level "Extension" row down to the "Extension.extension"
#   row's Constraint(s) column.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/based-on-value"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/based-on-value"

$constraintText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$elements.Range("AI2").Value = ""
$elements.Range("AI4").Value = $constraintText
